$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{ Row = 2; D = "306.92"; E = "-4.40%"; G = "4" },
    @{ Row = 3; D = "39.92"; E = "-5.92%"; G = "4" },
    @{ Row = 4; D = "5.019"; E = "-4.67%"; G = "4" },
    @{ Row = 5; D = "0.07670"; E = "-6.04%"; G = "4" },
    @{ Row = 6; D = "4.245"; E = "-2.32%"; G = "4" },
    @{ Row = 7; D = "1.619"; E = "-9.76%"; G = "4" },
    @{ Row = 8; D = "0.8849"; E = "-6.85%"; G = "4" },
    @{ Row = 9; D = "0.1002"; E = "-9.92%"; G = "4" },
    @{ Row = 10; D = "0.1724"; E = "-6.77%"; G = "4" },
    @{ Row = 11; D = "0.08931"; E = "-4.83%"; G = "4" },
    @{ Row = 12; D = "0.04389"; E = "-4.91%"; G = "4" },
    @{ Row = 13; E = "-0.52%"; G = "4" },
    @{ Row = 14; D = "0.001276"; E = "-0.58%"; G = "4" },
    @{ Row = 15; D = "0.005803"; E = "-1.60%"; G = "4" },
    @{ Row = 16; D = "3.352"; E = "-0.48%"; G = "4" },
    @{ Row = 17; D = "2.505"; E = "-0.22%"; G = "4" },
    @{ Row = 18; E = "-0.08%"; G = "4" },
    @{ Row = 19; D = "6.993"; E = "-5.93%"; G = "4" },
    @{ Row = 20; D = "0.1343"; E = "-3.77%"; G = "4" },
    @{ Row = 21; D = "0.3308"; E = "24.26%"; G = "4" },
    @{ Row = 22; D = "0.04214"; E = "0.58%"; G = "4" },
    @{ Row = 23; D = "0.001201"; E = "-4.49%"; G = "4" },
    @{ Row = 24; D = "0.004059"; E = "-5.34%"; G = "4" },
    @{ Row = 25; D = "0.0001223"; E = "-6.34%"; G = "4" },
    @{ Row = 26; E = "-0.37%"; G = "4" },
    @{ Row = 27; G = "4" },
    @{ Row = 28; G = "4" },
    @{ Row = 29; G = "4" },
    @{ Row = 30; G = "4" },
    @{ Row = 31; G = "4" },
    @{ Row = 32; G = "4" },
    @{ Row = 33; G = "4" },
    @{ Row = 34; G = "4" },
    @{ Row = 35; G = "4" },
    @{ Row = 36; G = "4" },
    @{ Row = 37; G = "4" },
    @{ Row = 38; D = "0.02345"; E = "-9.78%"; G = "4" },
    @{ Row = 39; D = "0.05172"; E = "-5.50%"; G = "4" },
    @{ Row = 40; D = "0.007968"; E = "2.22%"; G = "4" },
    @{ Row = 41; D = "0.1321"; E = "-5.24%"; G = "4" },
    @{ Row = 42; D = "0.006703"; E = "1.43%"; G = "4" },
    @{ Row = 43; D = "0.001998"; E = "-5.92%"; G = "4" },
    @{ Row = 44; D = "0.007656"; E = "-9.73%"; G = "4" },
    @{ Row = 45; D = "0.3046"; E = "-11.27%"; G = "4" },
    @{ Row = 46; D = "0.00006585"; E = "-5.65%"; G = "4" },
    @{ Row = 47; D = "0.00000000752"; E = "-0.19%"; G = "4" },
    @{ Row = 48; D = "0.003385"; E = "-2.59%"; G = "4" },
    @{ Row = 49; E = "41.07%"; G = "4" },
    @{ Row = 50; D = "0.00002106"; E = "-0.19%"; G = "4" },
    @{ Row = 51; D = "0.0002005"; E = "-0.19%"; G = "4" }
)

foreach ($item in $changes) {
    $r = $item.Row
    if ($item.ContainsKey("D")) {
        $cell = $ws.Cells.Item($r, 4)
        $cell.NumberFormat = "@"
        $cell.Value = $item.D
    }
    if ($item.ContainsKey("E")) {
        $cell = $ws.Cells.Item($r, 5)
        $cell.NumberFormat = "@"
        $cell.Value = $item.E
    }
    if ($item.ContainsKey("G")) {
        $cell = $ws.Cells.Item($r, 7)
        $cell.NumberFormat = "@"
        $cell.Value = $item.G
    }
}
